$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1394.3889
$ws.Range("J17").Value = 1394.3889
$ws.Range("L17").Value = 4183.1667
$ws.Range("N17").Value = -4519.1667
$ws.Range("H33").Value = 91.888885
$ws.Range("J33").Value = 104.545456
$ws.Range("L33").Value = 104.545456
$ws.Range("N33").Value = -562.5454560000001
$ws.Range("H64").Value = 3076
$ws.Range("I64").Value = 2910.4
$ws.Range("J64").Value = 3214
$ws.Range("K64").Value = 2910.4
$ws.Range("L64").Value = 3214
$ws.Range("M64").Value = -2662.4
$ws.Range("N64").Value = -3710
$ws.Range("H67").Value = 3076
$ws.Range("I67").Value = 2910.4
$ws.Range("J67").Value = 3214
$ws.Range("K67").Value = 2910.4
$ws.Range("L67").Value = 3214
$ws.Range("M67").Value = -2052.4
$ws.Range("N67").Value = -4930
$ws.Range("H80").Value = 876.3200000000001
$ws.Range("J80").Value = 983.8889
$ws.Range("L80").Value = 2951.6667
$ws.Range("N80").Value = -4947.6667
$ws.Range("H83").Value = 876.3200000000001
$ws.Range("J83").Value = 983.8889
$ws.Range("L83").Value = 8855.000100000001
$ws.Range("N83").Value = -18839.0001
$ws.Range("H88").Value = 3571.4285
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 3750
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 3750
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -4562
$ws.Range("H91").Value = 3571.4285
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 3750
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 3750
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -6558
$ws.Range("H100").Value = 3776.25
$ws.Range("I100").Value = 3776.25
$ws.Range("K100").Value = 3776.25
$ws.Range("M100").Value = -3235.25
$ws.Range("H112").Value = 6388.3335
$ws.Range("J112").Value = 6388.3335
$ws.Range("L112").Value = 19165.0005
$ws.Range("N112").Value = -21381.0005
$ws.Range("H135").Value = 681.3125
$ws.Range("I135").Value = 607.9231
$ws.Range("J135").Value = 999.3333
$ws.Range("K135").Value = 5471.3079
$ws.Range("L135").Value = 8993.9997
$ws.Range("M135").Value = -2936.3079
$ws.Range("N135").Value = -14063.9997
$ws.Range("H138").Value = 2358.7234
$ws.Range("I138").Value = 2365.3704
$ws.Range("K138").Value = 7096.111199999999
$ws.Range("M138").Value = -1956.111199999999
$ws.Range("H141").Value = 3159.9312
$ws.Range("I141").Value = 2496.5908
$ws.Range("K141").Value = 7489.7724
$ws.Range("M141").Value = -2309.7724

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1440.25
$ws.Range("I2").Value = 1437
$ws.Range("J2").Value = 1450
$ws.Range("K2").Value = 1437
$ws.Range("L2").Value = 1450
$ws.Range("M2").Value = -1324
$ws.Range("N2").Value = -1676
$ws.Range("H32").Value = 2682.8796
$ws.Range("I32").Value = 1878.5942
$ws.Range("J32").Value = 6646.857
$ws.Range("K32").Value = 1878.5942
$ws.Range("L32").Value = 6646.857
$ws.Range("M32").Value = -1591.5942
$ws.Range("N32").Value = -7220.857
$ws.Range("H61").Value = 3088.2
$ws.Range("I61").Value = 1571.5454
$ws.Range("K61").Value = 1571.5454
$ws.Range("M61").Value = -1359.5454
$ws.Range("H74").Value = 1277.3182
$ws.Range("I74").Value = 1174.0667
$ws.Range("J74").Value = 1498.5714
$ws.Range("K74").Value = 1174.0667
$ws.Range("L74").Value = 1498.5714
$ws.Range("M74").Value = -300.0667000000001
$ws.Range("N74").Value = -3246.5714
$ws.Range("H77").Value = 1277.3182
$ws.Range("I77").Value = 1174.0667
$ws.Range("J77").Value = 1498.5714
$ws.Range("K77").Value = 5870.333500000001
$ws.Range("L77").Value = 7492.857
$ws.Range("M77").Value = -1502.333500000001
$ws.Range("N77").Value = -16228.857
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21802
$ws.Range("H110").Value = 1650.9286
$ws.Range("I110").Value = 1007.6923
$ws.Range("K110").Value = 1007.6923
$ws.Range("M110").Value = 1037.3077
$ws.Range("H116").Value = 1440.25
$ws.Range("I116").Value = 1437
$ws.Range("J116").Value = 1450
$ws.Range("K116").Value = 1437
$ws.Range("L116").Value = 1450
$ws.Range("M116").Value = 857
$ws.Range("N116").Value = -6038
$ws.Range("H122").Value = 575
$ws.Range("I122").Value = 575
$ws.Range("K122").Value = 1725
$ws.Range("M122").Value = 725
$ws.Range("H136").Value = 3088.2
$ws.Range("I136").Value = 1571.5454
$ws.Range("K136").Value = 4714.6362
$ws.Range("M136").Value = -2164.6362

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1440.25
$ws.Range("I3").Value = 1437
$ws.Range("J3").Value = 1450
$ws.Range("K3").Value = 1437
$ws.Range("L3").Value = 1450
$ws.Range("M3").Value = -1323
$ws.Range("N3").Value = -1678
$ws.Range("H22").Value = 233
$ws.Range("I22").Value = 189.5
$ws.Range("K22").Value = 189.5
$ws.Range("M22").Value = -16.5
$ws.Range("H99").Value = 1837.6923
$ws.Range("J99").Value = 2300
$ws.Range("L99").Value = 2300
$ws.Range("N99").Value = -5296
$ws.Range("H105").Value = 2499.9092
$ws.Range("I105").Value = 2499.9092
$ws.Range("K105").Value = 2499.9092
$ws.Range("M105").Value = -752.9092000000001
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1936.3158
$ws.Range("I31").Value = 1614.9231
$ws.Range("J31").Value = 2632.6667
$ws.Range("K31").Value = 1614.9231
$ws.Range("L31").Value = 2632.6667
$ws.Range("M31").Value = -1319.9231
$ws.Range("N31").Value = -3222.6667
$ws.Range("H34").Value = 1936.3158
$ws.Range("I34").Value = 1614.9231
$ws.Range("J34").Value = 2632.6667
$ws.Range("K34").Value = 1614.9231
$ws.Range("L34").Value = 2632.6667
$ws.Range("M34").Value = -1412.9231
$ws.Range("N34").Value = -3036.6667
$ws.Range("H60").Value = 11747.25
$ws.Range("J60").Value = 11747.25
$ws.Range("L60").Value = 11747.25
$ws.Range("N60").Value = -12769.25
$ws.Range("H132").Value = 3067.2222
$ws.Range("I132").Value = 1759.6
$ws.Range("K132").Value = 5278.799999999999
$ws.Range("M132").Value = -2748.799999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 169.9
$ws.Range("I2").Value = 157
$ws.Range("K2").Value = 942
$ws.Range("M2").Value = -829
$ws.Range("H131").Value = 1380.58
$ws.Range("I131").Value = 593.8570999999999
$ws.Range("J131").Value = 1439.7957
$ws.Range("K131").Value = 1781.5713
$ws.Range("L131").Value = 4319.3871
$ws.Range("M131").Value = 3258.4287
$ws.Range("N131").Value = -14399.3871
$ws.Range("H140").Value = 1881.8214
$ws.Range("I140").Value = 1070.6875
$ws.Range("J140").Value = 2963.3333
$ws.Range("K140").Value = 3212.0625
$ws.Range("L140").Value = 8889.999899999999
$ws.Range("M140").Value = 1967.9375
$ws.Range("N140").Value = -19249.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 67
$ws.Range("I2").Value = 17.6
$ws.Range("K2").Value = 17.6
$ws.Range("M2").Value = 95.40000000000001
$ws.Range("H70").Value = 5350
$ws.Range("I70").Value = 7000
$ws.Range("J70").Value = 4800
$ws.Range("K70").Value = 7000
$ws.Range("L70").Value = 4800
$ws.Range("M70").Value = -6730
$ws.Range("N70").Value = -5340
$ws.Range("H73").Value = 5350
$ws.Range("I73").Value = 7000
$ws.Range("J73").Value = 4800
$ws.Range("K73").Value = 7000
$ws.Range("L73").Value = 4800
$ws.Range("M73").Value = -6064
$ws.Range("N73").Value = -6672
$ws.Range("H113").Value = 1235.1666
$ws.Range("I113").Value = 799.3333
$ws.Range("K113").Value = 799.3333
$ws.Range("M113").Value = 1370.6667
$ws.Range("H122").Value = 1808.1666
$ws.Range("I122").Value = 1462.4166
$ws.Range("J122").Value = 2499.6667
$ws.Range("K122").Value = 4387.2498
$ws.Range("L122").Value = 7499.000100000001
$ws.Range("M122").Value = -1937.2498
$ws.Range("N122").Value = -12399.0001
$ws.Range("H132").Value = 4518.077
$ws.Range("I132").Value = 4817.5713
$ws.Range("J132").Value = 4168.6665
$ws.Range("K132").Value = 14452.7139
$ws.Range("L132").Value = 12505.9995
$ws.Range("M132").Value = -11922.7139
$ws.Range("N132").Value = -17565.9995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 733
$ws.Range("I93").Value = 700
$ws.Range("K93").Value = 700
$ws.Range("M93").Value = 548
$ws.Range("H132").Value = 2733.3845
$ws.Range("I132").Value = 2139.6365
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 6418.9095
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -3888.9095
$ws.Range("N132").Value = -23057
$ws.Range("H136").Value = 4065.75
$ws.Range("I136").Value = 2479.8
$ws.Range("K136").Value = 7439.400000000001
$ws.Range("M136").Value = -4889.400000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 562.1429000000001
$ws.Range("I100").Value = 433.63635
$ws.Range("K100").Value = 867.2727
$ws.Range("M100").Value = -326.2727
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 9000
$ws.Range("M113").Value = -6830
$ws.Range("H122").Value = 53023.535
$ws.Range("I122").Value = 60951.848
$ws.Range("J122").Value = 1489.5
$ws.Range("K122").Value = 182855.544
$ws.Range("L122").Value = 4468.5
$ws.Range("M122").Value = -180405.544
$ws.Range("N122").Value = -9368.5
